# Ajout de la colonne "etat commande TMS" avec tracking number / statut "valide"
# (generation et ajout tracking number - elimine la commande deja creee en mode manuel)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nouvelle colonne Q : en-tete + valeur de la premiere commande
$ws.Range("Q1").Value = "etat commande TMS"
$ws.Range("Q2").Value = "valide"

# Largeur de la nouvelle colonne (au plus proche de ce que permet l'arrondi pixel du moteur)
$ws.Columns.Item(17).ColumnWidth = 19.83

# La nouvelle colonne devient la selection active de la feuille
$ws.Range("Q1:Q2").Select() | Out-Null
